$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E5").Value = "❌ EXPIRED 3760 days ago"
$ws.Range("E6").Value = "⚠️ Expires in 19 days"
$ws.Range("E7").Value = "⚠️ Expires in 19 days"
$ws.Range("E8").Value = "⚠️ Expires in 19 days"
